# Updates the "Pais" (paises.xlsx) COVID country stats sheet:
#  - refreshes the "Datos actualizados ..." timestamp in A1
#  - re-sorts several adjacent country pairs (names swap rows) and
#    refreshes the Casos/Recuperados/Muertes columns (B:H) with the
#    latest scrape for every affected row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 11:30"

$ws.Range("B7").Value = 368557
$ws.Range("C7").Value = 1293
$ws.Range("D7").Value = 194843
$ws.Range("E7").Value = 161440
$ws.Range("G7").Value = 12
$ws.Range("H7").Value = 12274

# Indonesia / Singapur swap rows (33<->34) + refreshed totals
$ws.Range("A33").Value = "Singapur"
$ws.Range("B33").Value = 41473
$ws.Range("C33").Value = 257
$ws.Range("D33").Value = 31938
$ws.Range("E33").Value = 9509
$ws.Range("H33").Value = 26

$ws.Range("A34").Value = "Indonesia"
$ws.Range("B34").Value = 41431
$ws.Range("D34").Value = 16243
$ws.Range("E34").Value = 22912
$ws.Range("H34").Value = 2276

$ws.Range("B40").Value = 31015
$ws.Range("C40").Value = 314
$ws.Range("D40").Value = 15317
$ws.Range("E40").Value = 14382
$ws.Range("G40").Value = 30
$ws.Range("H40").Value = 1316

# Afganistan / Filipinas swap rows (41<->42) + refreshed totals
$ws.Range("A41").Value = "Filipinas"
$ws.Range("B41").Value = 27799
$ws.Range("C41").Value = 561
$ws.Range("D41").Value = 7090
$ws.Range("E41").Value = 19593
$ws.Range("G41").Value = 8
$ws.Range("H41").Value = 1116

$ws.Range("A42").Value = "Afganistan"
$ws.Range("B42").Value = 27532
$ws.Range("C42").Value = 658
$ws.Range("D42").Value = 7660
$ws.Range("E42").Value = 19326
$ws.Range("G42").Value = 42
$ws.Range("H42").Value = 546

$ws.Range("B43").Value = 26818
$ws.Range("C43").Value = 739
$ws.Range("D43").Value = 13264
$ws.Range("E43").Value = 13435
$ws.Range("G43").Value = 3
$ws.Range("H43").Value = 119

$ws.Range("D66").Value = 7401
$ws.Range("E66").Value = 2442

$ws.Range("B70").Value = 8529
$ws.Range("C70").Value = 14
$ws.Range("D70").Value = 8000
$ws.Range("E70").Value = 408

# Tayikistan / Consejo Danes para los Refugiados swap rows (78<->79) + refreshed totals
$ws.Range("A78").Value = "Consejo Danes para los Refugiados"
$ws.Range("B78").Value = 5283
$ws.Range("C78").Value = 183
$ws.Range("D78").Value = 685
$ws.Range("E78").Value = 4481
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 117

$ws.Range("A79").Value = "Tayikistan"
$ws.Range("B79").Value = 5221
$ws.Range("D79").Value = 3700
$ws.Range("E79").Value = 1470
$ws.Range("H79").Value = 51

$ws.Range("D104").Value = 1421
$ws.Range("E104").Value = 492

# Lituania / Albania swap rows (110<->111) + refreshed totals
$ws.Range("A110").Value = "Albania"
$ws.Range("B110").Value = 1788
$ws.Range("C110").Value = 66
$ws.Range("D110").Value = 1086
$ws.Range("E110").Value = 663
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 39

$ws.Range("A111").Value = "Lituania"
$ws.Range("B111").Value = 1784
$ws.Range("C111").Value = 6
$ws.Range("D111").Value = 1449
$ws.Range("E111").Value = 259
$ws.Range("H111").Value = 76

$ws.Range("B113").Value = 1562
$ws.Range("C113").Value = 1
$ws.Range("D113").Value = 1443
$ws.Range("E113").Value = 91

# Nueva Zelanda / Eslovenia swap rows (115<->116) + refreshed totals
$ws.Range("A115").Value = "Eslovenia"
$ws.Range("B115").Value = 1511
$ws.Range("C115").Value = 8
$ws.Range("D115").Value = 1359
$ws.Range("E115").Value = 43
$ws.Range("H115").Value = 109

$ws.Range("A116").Value = "Nueva Zelanda"
$ws.Range("B116").Value = 1507
$ws.Range("C116").Value = 1
$ws.Range("D116").Value = 1482
$ws.Range("E116").Value = 3
$ws.Range("H116").Value = 22

$ws.Range("B125").Value = 1108
$ws.Range("C125").Value = 4
$ws.Range("D125").Value = 903
$ws.Range("E125").Value = 175

$ws.Range("B137").Value = 741
$ws.Range("C137").Value = 9
$ws.Range("D137").Value = 442
$ws.Range("E137").Value = 299

# Islas Feroe / Siria swap rows (166<->167) + refreshed totals
$ws.Range("A166").Value = "Siria"
$ws.Range("C166").Value = 9
$ws.Range("D166").Value = 78
$ws.Range("E166").Value = 102
$ws.Range("H166").Value = 7

$ws.Range("A167").Value = "Islas Feroe"
$ws.Range("B167").Value = 187
$ws.Range("D167").Value = 187
$ws.Range("E167").Value = 0
$ws.Range("H167").Value = 0

# Dominica / Fiyi swap rows (202<->203), totals unchanged
$ws.Range("A202").Value = "Fiyi"

$ws.Range("A203").Value = "Dominica"

# Santa Sede / Islas Turcas y Caicos swap rows (208<->209) + refreshed totals
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
